$d = $word.ActiveDocument

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function Insert-ParagraphXml([string]$innerXml) {
    $endPos = $d.Content.End
    $r = $d.Range($endPos, $endPos)
    $xml = '<w:p xmlns:w="' + $wNs + '">' + $innerXml + '</w:p>'
    $r.InsertXML($xml) | Out-Null
}

function Escape-Xml([string]$text) {
    $text = $text.Replace("&", "&amp;")
    $text = $text.Replace("<", "&lt;")
    $text = $text.Replace(">", "&gt;")
    return $text
}

function Insert-PlainParagraph([string]$text) {
    $t = Escape-Xml $text
    Insert-ParagraphXml ("<w:r><w:t>" + $t + "</w:t></w:r>")
}

function Insert-ListParagraph([string]$text) {
    $t = Escape-Xml $text
    $pPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr>'
    Insert-ParagraphXml ($pPr + "<w:r><w:t>" + $t + "</w:t></w:r>")
}

Insert-PlainParagraph "15/5"
Insert-ListParagraph "Added acronyms and fixed custom captions for list of figures"
Insert-PlainParagraph "16/5"
Insert-ListParagraph "Poster and Presentation first draft"
Insert-ParagraphXml ""
